$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table name / description (row 1-2) ---
$ws.Range("B1").Value = "Movimentacoes"

# --- Rows 5-9: rename columns / types ---
$ws.Range("A5").Value = "id"
$ws.Range("A6").Value = "vaga_id"
$ws.Range("H6").Value = "Foreign Key da tabela vagas"
$ws.Range("A7").Value = "veiculo_id"
$ws.Range("H7").Value = "Foreign Key da tabela veiculos"
$ws.Range("A8").Value = "entrada"
$ws.Range("C8").Value = "timestamp"
$ws.Range("A9").Value = "saida"
$ws.Range("C9").Value = "timestamp"

# --- Row 10: clear the valor_pago row contents ---
$ws.Range("A10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("H10").Value = ""

# --- Rows 13-15: add index information ---
$ws.Range("A13").Value = "PRIMARY"
$ws.Range("C13").Value = "Sim"
$ws.Range("D13").Value = "Não"
$ws.Range("E13").Value = "Sim"
$ws.Range("F13").Value = "id"

$ws.Range("A14").Value = "Index_vaga_id"
$ws.Range("C14").Value = "Não"
$ws.Range("D14").Value = "Sim"
$ws.Range("E14").Value = "Não"
$ws.Range("F14").Value = "vaga_id"

$ws.Range("A15").Value = "index_veiculo_id"
$ws.Range("C15").Value = "Não"
$ws.Range("D15").Value = "Sim"
$ws.Range("E15").Value = "Não"
$ws.Range("F15").Value = "veiculo_id"

# Normalize row 15 styles (F/G/H) to match rows 13/14/16
$ws.Range("F14:H14").Copy()
$ws.Range("F15:H15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Selection change ---
$ws.Range("A16:B16").Select()

$wb.Save()
